# Fruta / hortaliza, semanal
# Insert two new weekly-report rows for "Chirimoya" (Cultivar IV Región) at
# rows 580-581, pushing the existing rows (old 580..610) down to 582..612.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 580 (each Insert() pushes
# everything at/below row 580 down by one, so calling it twice opens a
# two-row gap at 580-581 and the old data that lived there now starts at 582).
$ws.Rows.Item(580).Insert()
$ws.Rows.Item(580).Insert()

# New row 580: Chirimoya, Cultivar IV Región, calidad "Primera"
$ws.Range("A580").Value = 6
$ws.Range("B580").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C580").Value = "Metropolitana"
$ws.Range("D580").Value = 45267
$ws.Range("E580").Value = 13
$ws.Range("F580").Value = "Fruta"
$ws.Range("G580").Value = 100107
$ws.Range("H580").Value = "Otros"
$ws.Range("I580").Value = 100107002
$ws.Range("J580").Value = "Chirimoya"
$ws.Range("K580").Value = "Cultivar IV Región"
$ws.Range("L580").Value = "Primera"
$ws.Range("M580").Value = 970
$ws.Range("N580").Value = 16000
$ws.Range("O580").Value = 16000
$ws.Range("P580").Value = 16000
$ws.Range("Q580").Value = "`$/bandeja 10 kilos"
$ws.Range("R580").Value = "Provincia de Limarí"
$ws.Range("S580").Value = 1600
$ws.Range("T580").Value = 10

# New row 581: Chirimoya, Cultivar IV Región, calidad "Segunda"
$ws.Range("A581").Value = 6
$ws.Range("B581").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C581").Value = "Metropolitana"
$ws.Range("D581").Value = 45267
$ws.Range("E581").Value = 13
$ws.Range("F581").Value = "Fruta"
$ws.Range("G581").Value = 100107
$ws.Range("H581").Value = "Otros"
$ws.Range("I581").Value = 100107002
$ws.Range("J581").Value = "Chirimoya"
$ws.Range("K581").Value = "Cultivar IV Región"
$ws.Range("L581").Value = "Segunda"
$ws.Range("M581").Value = 800
$ws.Range("N581").Value = 13000
$ws.Range("O581").Value = 13000
$ws.Range("P581").Value = 13000
$ws.Range("Q581").Value = "`$/bandeja 10 kilos"
$ws.Range("R581").Value = "Provincia de Limarí"
$ws.Range("S581").Value = 1300
$ws.Range("T581").Value = 10
